$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7819745094127484
$ws.Range("C2").Value = 0.2335998566035755
$ws.Range("E2").Value = 0.4264233615170383
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002364026894831113
$ws.Range("I2").Value = 0.2574922939564068
$ws.Range("N2").Value = 0.7925338793933321
$ws.Range("O2").Value = 1.23301000020362
$ws.Range("B3").Value = 0.6834810727099807
$ws.Range("C3").Value = 0.2059933084904344
$ws.Range("E3").Value = 0.3719117739346842
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.00236688649945028
$ws.Range("I3").Value = 0.2613142563083848
$ws.Range("N3").Value = 0.7927348538646513
$ws.Range("O3").Value = 1.220348518214905
$ws.Range("B4").Value = 0.6228653617803843
$ws.Range("C4").Value = 0.1889585801924625
$ws.Range("E4").Value = 0.338534222790301
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002368733479126359
$ws.Range("I4").Value = 0.2639244056750307
$ws.Range("N4").Value = 0.7932216771413536
$ws.Range("O4").Value = 1.213853809957669
$ws.Range("B5").Value = 0.598129734498599
$ws.Range("C5").Value = 0.1819957606428488
$ws.Range("E5").Value = 0.3249540513085094
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002369509138759968
$ws.Range("I5").Value = 0.2650540123980285
$ws.Range("N5").Value = 0.7935116613955628
$ws.Range("O5").Value = 1.211527034112919
$ws.Range("B6").Value = 0.5940203661940018
$ws.Range("C6").Value = 0.1808383239067837
$ws.Range("E6").Value = 0.3227003151194339
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002369639327621988
$ws.Range("I6").Value = 0.2652455572543637
$ws.Range("N6").Value = 0.7935653518793018
$ws.Range("O6").Value = 1.211159939588555
$ws.Range("B7").Value = 0.6225319056602814
$ws.Range("C7").Value = 0.1888647621275652
$ws.Range("E7").Value = 0.3383509914541918
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002368743846688401
$ws.Range("I7").Value = 0.2639393733136224
$ws.Range("N7").Value = 0.7932252168016518
$ws.Range("O7").Value = 1.213821137583096
$ws.Range("B8").Value = 0.7480437027278981
$ws.Range("C8").Value = 0.2240986703062902
$ws.Range("E8").Value = 0.4076073590350404
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.00236499400907926
$ws.Range("I8").Value = 0.2587552491693152
$ws.Range("N8").Value = 0.7925278709727621
$ws.Range("O8").Value = 1.228377628650321
$ws.Range("B9").Value = 0.9930240884717705
$ws.Range("C9").Value = 0.2925213008762739
$ws.Range("E9").Value = 0.5442463263346724
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002358360569978073
$ws.Range("I9").Value = 0.2506929228981249
$ws.Range("N9").Value = 0.7940357685497617
$ws.Range("O9").Value = 1.267163016626455
$ws.Range("B10").Value = 1.172283230444407
$ws.Range("C10").Value = 0.3423834500133296
$ws.Range("E10").Value = 0.6452750528707583
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002353921097370253
$ws.Range("I10").Value = 0.2460705199229167
$ws.Range("N10").Value = 0.7968872246758423
$ws.Range("O10").Value = 1.302028356160918
$ws.Range("B11").Value = 1.253670051714323
$ws.Range("C11").Value = 0.3649790310422532
$ws.Range("E11").Value = 0.6914046423499656
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002351994704321794
$ws.Range("I11").Value = 0.2442541671303395
$ws.Range("N11").Value = 0.798561290646802
$ws.Range("O11").Value = 1.319300413951254
$ws.Range("B12").Value = 1.284465461197556
$ws.Range("C12").Value = 0.3735228245896565
$ws.Range("E12").Value = 0.7088995885753917
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.00235127854600603
$ws.Range("I12").Value = 0.2436078760519287
$ws.Range("N12").Value = 0.799249259824748
$ws.Range("O12").Value = 1.326046049117338
$ws.Range("B13").Value = 1.277834204816429
$ws.Range("C13").Value = 0.3716833314300629
$ws.Range("E13").Value = 0.7051305195087139
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002351432191935703
$ws.Range("I13").Value = 0.243745214285326
$ws.Range("N13").Value = 0.799098692598804
$ws.Range("O13").Value = 1.324584101439001
$ws.Range("B14").Value = 1.25620409634729
$ws.Range("C14").Value = 0.3656821887480817
$ws.Range("E14").Value = 0.6928434156457826
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002351935518840609
$ws.Range("I14").Value = 0.2442001622643168
$ws.Range("N14").Value = 0.7986168082263134
$ws.Range("O14").Value = 1.319851259268461
$ws.Range("B15").Value = 1.242951869426122
$ws.Range("C15").Value = 0.3620046619080597
$ws.Range("E15").Value = 0.6853207449128149
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002352245554495849
$ws.Range("I15").Value = 0.2444842485240564
$ws.Range("N15").Value = 0.798328672547953
$ws.Range("O15").Value = 1.3169790266642
$ws.Range("B16").Value = 1.166961111810792
$ws.Range("C16").Value = 0.3409050122208157
$ws.Range("E16").Value = 0.6422639950437343
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002354048859971948
$ws.Range("I16").Value = 0.2461950148539565
$ws.Range("N16").Value = 0.7967853909817819
$ws.Range("O16").Value = 1.300928170767975
$ws.Range("B17").Value = 1.120301689076143
$ws.Range("C17").Value = 0.327938674526564
$ws.Range("E17").Value = 0.6158953628066257
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002355178936536806
$ws.Range("I17").Value = 0.2473181090596803
$ws.Range("N17").Value = 0.7959350577736188
$ws.Range("O17").Value = 1.291444529108247
$ws.Range("B18").Value = 1.093449519120838
$ws.Range("C18").Value = 0.3204726041081187
$ws.Range("E18").Value = 0.6007448043084196
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002355837698291098
$ws.Range("I18").Value = 0.2479910203075804
$ws.Range("N18").Value = 0.7954814558691936
$ws.Range("O18").Value = 1.286122533925635
$ws.Range("B19").Value = 1.084355310887929
$ws.Range("C19").Value = 0.3179433215864265
$ws.Range("E19").Value = 0.5956177677434624
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002356062252384202
$ws.Range("I19").Value = 0.2482234726994932
$ws.Range("N19").Value = 0.7953339741644498
$ws.Range("O19").Value = 1.28434333333098
$ws.Range("B20").Value = 1.125270215093451
$ws.Range("C20").Value = 0.3293198106417776
$ws.Range("E20").Value = 0.6187006798457872
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002355057730839899
$ws.Range("I20").Value = 0.2471957634211712
$ws.Range("N20").Value = 0.7960219056065228
$ws.Range("O20").Value = 1.292440325263101
$ws.Range("B21").Value = 1.262558044281263
$ws.Range("C21").Value = 0.3674452147227498
$ws.Range("E21").Value = 0.6964516954516995
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002351787318582521
$ws.Range("I21").Value = 0.2440654034101968
$ws.Range("N21").Value = 0.7987568840822377
$ws.Range("O21").Value = 1.321235827176224
$ws.Range("B22").Value = 1.352143086587375
$ws.Range("C22").Value = 0.3922884921168475
$ws.Range("E22").Value = 0.7474229179322123
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.00234972754942823
$ws.Range("I22").Value = 0.2422617137831153
$ws.Range("N22").Value = 0.8008592466483719
$ws.Range("O22").Value = 1.341251702567604
$ws.Range("B23").Value = 1.30434309254872
$ws.Range("C23").Value = 0.3790359785195392
$ws.Range("E23").Value = 0.7202036143974624
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002350819806255345
$ws.Range("I23").Value = 0.243202105200016
$ws.Range("N23").Value = 0.7997084157155854
$ws.Range("O23").Value = 1.330458693986429
$ws.Range("B24").Value = 1.123024028644352
$ws.Range("C24").Value = 0.3286954349750033
$ws.Range("E24").Value = 0.6174323678046534
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002355112499620341
$ws.Range("I24").Value = 0.2472509910880092
$ws.Range("N24").Value = 0.7959825318303899
$ws.Range("O24").Value = 1.29198972015962
$ws.Range("B25").Value = 0.9268762410067666
$ws.Range("C25").Value = 0.2740828974629039
$ws.Range("E25").Value = 0.507179335981661
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002360078510230196
$ws.Range("I25").Value = 0.2526468698234403
$ws.Range("N25").Value = 0.7933210022432888
$ws.Range("O25").Value = 1.255561165833257

Write-Output "Updated 192 cells (B,C,E,F,G,I,N,O) for rows 2-25 for case with 380 kV"
